# no-op test
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
